# Updated cryptos list on Wed May 24 16:40:23 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row (rows 2-51) on the active worksheet with the latest
# scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price, Volume(1h)). A $null Price means that row's Price column
# is unchanged - only the Volume(1h) value moved for that coin.
$updates = @{
    2  = @("26.398.68",     "  -3.47%  ")
    3  = @("1.801.22",      "  -3.10%  ")
    4  = @("1.008",         "  +0.56%  ")
    5  = @($null,           "  +0.47%  ")
    6  = @("307.81",        "  -2.27%  ")
    7  = @("0.4504",        "  -2.38%  ")
    8  = @("0.3646",        "  -1.89%  ")
    9  = @("0.07088",       "  -3.01%  ")
    10 = @("0.8662",        "  -2.91%  ")
    11 = @("0.07770",       "  -0.77%  ")
    12 = @("19.23",         "  -4.37%  ")
    13 = @("1.834.06",      "  +2.33%  ")
    14 = @($null,           "  -2.86%  ")
    15 = @("6.296",         "  -3.79%  ")
    16 = @("86.13",         "  -6.05%  ")
    17 = @("1.009",         "  +0.48%  ")
    18 = @("0.000008527",   "  -4.57%  ")
    19 = @($null,           "  +0.45%  ")
    20 = @("26.442.62",     "  -3.36%  ")
    21 = @("14.19",         "  -4.08%  ")
    22 = @("4.946",         "  -3.58%  ")
    23 = @("10.39",         "  -1.73%  ")
    24 = @($null,           "  +2.34%  ")
    25 = @("149.37",        "  -1.92%  ")
    26 = @("17.86",         "  -3.39%  ")
    27 = @("1.969",         "  -4.44%  ")
    28 = @("112.58",        "  -3.12%  ")
    29 = @("4.847",         "  -4.72%  ")
    30 = @("0.08620",       "  -2.40%  ")
    31 = @("3.025",         "  -1.92%  ")
    32 = @("0.7266",        "  -6.15%  ")
    33 = @("4.423",         "  -2.09%  ")
    34 = @("1.110",         "  -5.62%  ")
    35 = @("2.536",         "  -6.81%  ")
    36 = @("1.071",         "  -0.80%  ")
    37 = @("0.01914",       "  -2.36%  ")
    38 = @("0.05054",       "  -4.18%  ")
    39 = @("2.872",         "  -3.22%  ")
    40 = @("6.951",         "  -1.31%  ")
    41 = @("0.4891",        "  -4.87%  ")
    42 = @("0.1562",        "  -4.96%  ")
    43 = @("8.092",         "  -4.29%  ")
    44 = @($null,           "  +0.54%  ")
    45 = @("0.4597",        "  -4.40%  ")
    46 = @("101.04",        "  -1.58%  ")
    47 = @("9.860",         "  -5.08%  ")
    48 = @("1.578",         "  -4.18%  ")
    49 = @("0.05983",       "  -3.83%  ")
    50 = @("63.23",         "  -3.85%  ")
    51 = @("36.05",         "  -2.01%  ")
}

# The Price column holds values like "1.008" or "0.4504" that Excel's COM
# layer would otherwise auto-convert to numbers (dropping the exact text
# representation, e.g. trailing zeros). Force the whole Price column to a
# text format first so every assigned value is stored verbatim as a
# string, then restore the default "Normal" style afterwards so the
# on-disk cell formatting is unchanged from the original workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($null -ne $price) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    $ws.Cells.Item($row, 5).Value = $volume
}

$priceRange.Style = "Normal"
